$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SAEX")
$ws.Columns("D").Insert()
$ws.Range("D7").Value = 43465
